# Auto-generated edit script applying the Leviathan_Profits market-data refresh
# (scheduled runner updates currentAveragePrice* / LevePrice* / LeveProfit* columns)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 512.9787  # ALC!H15
$ws.Cells.Item(15, 9).Value = 512.9787  # ALC!I15
$ws.Cells.Item(15, 11).Value = 1538.9361  # ALC!K15
$ws.Cells.Item(15, 13).Value = -1369.9361  # ALC!M15

$ws.Cells.Item(17, 8).Value = 456305.4  # ALC!H17
$ws.Cells.Item(17, 9).Value = 1620  # ALC!I17
$ws.Cells.Item(17, 10).Value = 477957.1  # ALC!J17
$ws.Cells.Item(17, 11).Value = 4860  # ALC!K17
$ws.Cells.Item(17, 12).Value = 1433871.3  # ALC!L17
$ws.Cells.Item(17, 13).Value = -4692  # ALC!M17
$ws.Cells.Item(17, 14).Value = -1434207.3  # ALC!N17

$ws.Cells.Item(33, 8).Value = 500  # ALC!H33
$ws.Cells.Item(33, 9).Value = 0  # ALC!I33
$ws.Cells.Item(33, 11).Value = 0  # ALC!K33
$ws.Cells.Item(33, 13).ClearContents()  # ALC!M33

$ws.Cells.Item(57, 8).Value = 34749.25  # ALC!H57
$ws.Cells.Item(57, 10).Value = 34749.25  # ALC!J57
$ws.Cells.Item(57, 12).Value = 104247.75  # ALC!L57
$ws.Cells.Item(57, 14).Value = -105245.75  # ALC!N57

$ws.Cells.Item(95, 8).Value = 38973.4  # ALC!H95
$ws.Cells.Item(95, 10).Value = 38973.4  # ALC!J95
$ws.Cells.Item(95, 12).Value = 38973.4  # ALC!L95
$ws.Cells.Item(95, 14).Value = -44465.4  # ALC!N95

$ws.Cells.Item(113, 8).Value = 5250  # ALC!H113
$ws.Cells.Item(113, 9).Value = 3666.6667  # ALC!I113
$ws.Cells.Item(113, 10).Value = 6200  # ALC!J113
$ws.Cells.Item(113, 11).Value = 3666.6667  # ALC!K113
$ws.Cells.Item(113, 12).Value = 6200  # ALC!L113
$ws.Cells.Item(113, 13).Value = -412.6667000000002  # ALC!M113
$ws.Cells.Item(113, 14).Value = -12708  # ALC!N113

$ws.Cells.Item(116, 8).Value = 4576  # ALC!H116
$ws.Cells.Item(116, 9).Value = 4289.8  # ALC!I116
$ws.Cells.Item(116, 11).Value = 4289.8  # ALC!K116
$ws.Cells.Item(116, 13).Value = -847.8000000000002  # ALC!M116

$ws.Cells.Item(135, 8).Value = 48120.715  # ALC!H135
$ws.Cells.Item(135, 9).Value = 525.4286  # ALC!I135
$ws.Cells.Item(135, 10).Value = 143311.28  # ALC!J135
$ws.Cells.Item(135, 11).Value = 4728.8574  # ALC!K135
$ws.Cells.Item(135, 12).Value = 1289801.52  # ALC!L135
$ws.Cells.Item(135, 13).Value = -2193.8574  # ALC!M135
$ws.Cells.Item(135, 14).Value = -1294871.52  # ALC!N135

$ws.Cells.Item(137, 8).Value = 1966.2174  # ALC!H137
$ws.Cells.Item(137, 9).Value = 1811.15  # ALC!I137
$ws.Cells.Item(137, 11).Value = 5433.450000000001  # ALC!K137
$ws.Cells.Item(137, 13).Value = -2883.450000000001  # ALC!M137

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 2211.0833  # ARM!H122
$ws.Cells.Item(122, 9).Value = 2146.4644  # ARM!I122
$ws.Cells.Item(122, 10).Value = 2437.25  # ARM!J122
$ws.Cells.Item(122, 11).Value = 6439.3932  # ARM!K122
$ws.Cells.Item(122, 12).Value = 7311.75  # ARM!L122
$ws.Cells.Item(122, 13).Value = -3989.3932  # ARM!M122
$ws.Cells.Item(122, 14).Value = -12211.75  # ARM!N122

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 8780.947  # BSM!H20
$ws.Cells.Item(20, 9).Value = 8426.083000000001  # BSM!I20
$ws.Cells.Item(20, 10).Value = 9389.286  # BSM!J20
$ws.Cells.Item(20, 11).Value = 8426.083000000001  # BSM!K20
$ws.Cells.Item(20, 12).Value = 9389.286  # BSM!L20
$ws.Cells.Item(20, 13).Value = -8179.083000000001  # BSM!M20
$ws.Cells.Item(20, 14).Value = -9883.286  # BSM!N20

$ws.Cells.Item(105, 8).Value = 1458.4  # BSM!H105
$ws.Cells.Item(105, 10).Value = 0  # BSM!J105
$ws.Cells.Item(105, 12).Value = 0  # BSM!L105
$ws.Cells.Item(105, 14).ClearContents()  # BSM!N105

$ws.Cells.Item(134, 8).Value = 885.5833  # BSM!H134
$ws.Cells.Item(134, 9).Value = 875.1818  # BSM!I134
$ws.Cells.Item(134, 11).Value = 2625.5454  # BSM!K134
$ws.Cells.Item(134, 13).Value = -90.54539999999997  # BSM!M134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 500001000  # CRP!H7
$ws.Cells.Item(7, 9).Value = 500001000  # CRP!I7
$ws.Cells.Item(7, 11).Value = 500001000  # CRP!K7
$ws.Cells.Item(7, 13).Value = -500000887  # CRP!M7

$ws.Cells.Item(58, 8).Value = 1862.762  # CRP!H58
$ws.Cells.Item(58, 9).Value = 1844.1111  # CRP!I58
$ws.Cells.Item(58, 10).Value = 1974.6666  # CRP!J58
$ws.Cells.Item(58, 11).Value = 1844.1111  # CRP!K58
$ws.Cells.Item(58, 12).Value = 1974.6666  # CRP!L58
$ws.Cells.Item(58, 13).Value = -1641.1111  # CRP!M58
$ws.Cells.Item(58, 14).Value = -2380.6666  # CRP!N58

$ws.Cells.Item(86, 8).Value = 4894.75  # CRP!H86
$ws.Cells.Item(86, 9).Value = 4821.5  # CRP!I86
$ws.Cells.Item(86, 10).Value = 4919.1665  # CRP!J86
$ws.Cells.Item(86, 11).Value = 4821.5  # CRP!K86
$ws.Cells.Item(86, 12).Value = 4919.1665  # CRP!L86
$ws.Cells.Item(86, 13).Value = -3698.5  # CRP!M86
$ws.Cells.Item(86, 14).Value = -7165.1665  # CRP!N86

$ws.Cells.Item(89, 8).Value = 4894.75  # CRP!H89
$ws.Cells.Item(89, 9).Value = 4821.5  # CRP!I89
$ws.Cells.Item(89, 10).Value = 4919.1665  # CRP!J89
$ws.Cells.Item(89, 11).Value = 24107.5  # CRP!K89
$ws.Cells.Item(89, 12).Value = 24595.8325  # CRP!L89
$ws.Cells.Item(89, 13).Value = -18491.5  # CRP!M89
$ws.Cells.Item(89, 14).Value = -35827.8325  # CRP!N89

$ws.Cells.Item(122, 8).Value = 13333  # CRP!H122
$ws.Cells.Item(122, 9).Value = 13333  # CRP!I122
$ws.Cells.Item(122, 11).Value = 39999  # CRP!K122
$ws.Cells.Item(122, 13).Value = -37549  # CRP!M122

$ws.Cells.Item(132, 8).Value = 2941.9285  # CRP!H132
$ws.Cells.Item(132, 9).Value = 2926.0908  # CRP!I132
$ws.Cells.Item(132, 11).Value = 8778.2724  # CRP!K132
$ws.Cells.Item(132, 13).Value = -6248.2724  # CRP!M132

$ws.Cells.Item(134, 8).Value = 3764.3215  # CRP!H134
$ws.Cells.Item(134, 9).Value = 3795.9565  # CRP!I134
$ws.Cells.Item(134, 11).Value = 11387.8695  # CRP!K134
$ws.Cells.Item(134, 13).Value = -8852.869499999999  # CRP!M134

$ws.Cells.Item(136, 8).Value = 1862.762  # CRP!H136
$ws.Cells.Item(136, 9).Value = 1844.1111  # CRP!I136
$ws.Cells.Item(136, 10).Value = 1974.6666  # CRP!J136
$ws.Cells.Item(136, 11).Value = 5532.3333  # CRP!K136
$ws.Cells.Item(136, 12).Value = 5923.9998  # CRP!L136
$ws.Cells.Item(136, 13).Value = -2982.3333  # CRP!M136
$ws.Cells.Item(136, 14).Value = -11023.9998  # CRP!N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(88, 8).Value = 10333.333  # CUL!H88
$ws.Cells.Item(88, 10).Value = 10333.333  # CUL!J88
$ws.Cells.Item(88, 12).Value = 30999.999  # CUL!L88
$ws.Cells.Item(88, 14).Value = -31855.999  # CUL!N88

$ws.Cells.Item(91, 8).Value = 10333.333  # CUL!H91
$ws.Cells.Item(91, 10).Value = 10333.333  # CUL!J91
$ws.Cells.Item(91, 12).Value = 30999.999  # CUL!L91
$ws.Cells.Item(91, 14).Value = -33963.999  # CUL!N91

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(47, 8).Value = 52999  # GSM!H47
$ws.Cells.Item(47, 10).Value = 52999  # GSM!J47
$ws.Cells.Item(47, 12).Value = 52999  # GSM!L47
$ws.Cells.Item(47, 14).Value = -54135  # GSM!N47

$ws.Cells.Item(122, 8).Value = 2049.3333  # GSM!H122
$ws.Cells.Item(122, 9).Value = 1725.1111  # GSM!I122
$ws.Cells.Item(122, 10).Value = 2535.6667  # GSM!J122
$ws.Cells.Item(122, 11).Value = 5175.3333  # GSM!K122
$ws.Cells.Item(122, 12).Value = 7607.000100000001  # GSM!L122
$ws.Cells.Item(122, 13).Value = -2725.3333  # GSM!M122
$ws.Cells.Item(122, 14).Value = -12507.0001  # GSM!N122

$ws.Cells.Item(132, 8).Value = 3541.0278  # GSM!H132
$ws.Cells.Item(132, 9).Value = 2326.2307  # GSM!I132
$ws.Cells.Item(132, 10).Value = 6699.5  # GSM!J132
$ws.Cells.Item(132, 11).Value = 6978.6921  # GSM!K132
$ws.Cells.Item(132, 12).Value = 20098.5  # GSM!L132
$ws.Cells.Item(132, 13).Value = -4448.6921  # GSM!M132
$ws.Cells.Item(132, 14).Value = -25158.5  # GSM!N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4416.914  # LTW!H40
$ws.Cells.Item(40, 9).Value = 3784.9614  # LTW!I40
$ws.Cells.Item(40, 11).Value = 3784.9614  # LTW!K40
$ws.Cells.Item(40, 13).Value = -3648.9614  # LTW!M40

$ws.Cells.Item(122, 8).Value = 10492.929  # LTW!H122
$ws.Cells.Item(122, 9).Value = 12440.2  # LTW!I122
$ws.Cells.Item(122, 11).Value = 37320.60000000001  # LTW!K122
$ws.Cells.Item(122, 13).Value = -34870.60000000001  # LTW!M122

$ws.Cells.Item(130, 8).Value = 74714.5  # LTW!H130
$ws.Cells.Item(130, 10).Value = 74714.5  # LTW!J130
$ws.Cells.Item(130, 12).Value = 74714.5  # LTW!L130
$ws.Cells.Item(130, 14).Value = -84754.5  # LTW!N130

$ws.Cells.Item(132, 8).Value = 4165.3887  # LTW!H132
$ws.Cells.Item(132, 9).Value = 3819.4  # LTW!I132
$ws.Cells.Item(132, 10).Value = 4298.4614  # LTW!J132
$ws.Cells.Item(132, 11).Value = 11458.2  # LTW!K132
$ws.Cells.Item(132, 12).Value = 12895.3842  # LTW!L132
$ws.Cells.Item(132, 13).Value = -8928.200000000001  # LTW!M132
$ws.Cells.Item(132, 14).Value = -17955.3842  # LTW!N132

$ws.Cells.Item(136, 8).Value = 4810.16  # LTW!H136
$ws.Cells.Item(136, 9).Value = 4263.7334  # LTW!I136
$ws.Cells.Item(136, 10).Value = 5629.8  # LTW!J136
$ws.Cells.Item(136, 11).Value = 12791.2002  # LTW!K136
$ws.Cells.Item(136, 12).Value = 16889.4  # LTW!L136
$ws.Cells.Item(136, 13).Value = -10241.2002  # LTW!M136
$ws.Cells.Item(136, 14).Value = -21989.4  # LTW!N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 10644.8  # WVR!H4
$ws.Cells.Item(4, 9).Value = 24581.666  # WVR!I4
$ws.Cells.Item(4, 10).Value = 4671.857  # WVR!J4
$ws.Cells.Item(4, 11).Value = 24581.666  # WVR!K4
$ws.Cells.Item(4, 12).Value = 4671.857  # WVR!L4
$ws.Cells.Item(4, 13).Value = -24468.666  # WVR!M4
$ws.Cells.Item(4, 14).Value = -4897.857  # WVR!N4

$ws.Cells.Item(122, 8).Value = 1999.5  # WVR!H122
$ws.Cells.Item(122, 9).Value = 1999.5  # WVR!I122
$ws.Cells.Item(122, 11).Value = 5998.5  # WVR!K122
$ws.Cells.Item(122, 13).Value = -3548.5  # WVR!M122

$ws.Cells.Item(123, 8).Value = 0  # WVR!H123
$ws.Cells.Item(123, 10).Value = 0  # WVR!J123
$ws.Cells.Item(123, 12).Value = 0  # WVR!L123
$ws.Cells.Item(123, 14).ClearContents()  # WVR!N123

$ws.Cells.Item(132, 8).Value = 4630.0884  # WVR!H132
$ws.Cells.Item(132, 9).Value = 5278.8887  # WVR!I132
$ws.Cells.Item(132, 10).Value = 2127.5715  # WVR!J132
$ws.Cells.Item(132, 11).Value = 15836.6661  # WVR!K132
$ws.Cells.Item(132, 12).Value = 6382.7145  # WVR!L132
$ws.Cells.Item(132, 13).Value = -13306.6661  # WVR!M132
$ws.Cells.Item(132, 14).Value = -11442.7145  # WVR!N132
